# Admin page update: split the single "Email" column layout into
# Timestamp / Date / Time / Email, and refresh the sample log row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 holds a date-looking literal ("3/7/2025") that must stay as TEXT,
# not get auto-converted into a date serial number. Mark it as Text
# before assigning, then flip the cell back to the Normal style so no
# stray number format lingers on the finished cell.
$ws.Range("B2").NumberFormat = "@"

$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Time"
$ws.Range("D1").Value = "Email"

$ws.Range("A2").Value = "2025-07-03T08:40:04.963Z"
$ws.Range("B2").Value = "3/7/2025"
$ws.Range("C2").Value = "1:40:04 am"
$ws.Range("D2").Value = "kunaldutt69@gmail.com"

$ws.Range("B2").Style = "Normal"
